$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# The paragraph that talks about "fast-forward errors" ends with a
# bookmarked run pair:  "then be" + " able to push their changes to the
# remote repository."   We need to turn that into three runs:
#   "then "  /  "pull the remote repository and after that they should
#   be "  /  "able to push their changes to the remote repository."
# and drop the (hidden) "_GoBack" bookmark that wraps the split point.
# ----------------------------------------------------------------------

$oldMiddle = "then be"
$oldTail   = " able to push their changes to the remote repository."

$newRun1 = "then "
$newRun2 = "pull the remote repository and after that they should be "
$newRun3 = "able to push their changes to the remote repository."

$full = $d.Content.Text
$idx = $full.IndexOf($oldMiddle)
if ($idx -lt 0) {
    throw "Could not locate '$oldMiddle' in the document"
}

$run2End = $idx + $oldMiddle.Length
$run3Len = $oldTail.Length

# sanity check: the text right after "then be" must be the expected tail
$tailCheck = $d.Range($run2End, $run2End + $run3Len).Text
if ($tailCheck -ne $oldTail) {
    throw "Unexpected tail text: [$tailCheck]"
}

# 1) Bookmark the run2/run3 boundary so subsequent edits can't bleed
#    across it and silently re-merge the runs we are about to create.
$d.Bookmarks.Add("ZZ_SPLIT_M", $d.Range($run2End, $run2End)) | Out-Null

# 2) Replace run2 ("then be") in full -> becomes a clean "then " run.
$rng2 = $d.Range($idx, $run2End)
$rng2.Text = $newRun1

# 3) Replace run3 (" able to push...") in full with the combined text
#    for the two new trailing runs -> another clean run.
$bmM = $d.Bookmarks("ZZ_SPLIT_M")
$rng3 = $d.Range($bmM.Start, $bmM.Start + $run3Len)
$rng3.Text = $newRun2 + $newRun3

# 4) Split that combined run into its two final pieces the same way:
#    bookmark the split point, then do two full-range replacements.
$bmM2 = $d.Bookmarks("ZZ_SPLIT_M")
$run3Start = $bmM2.Start
$run3EndNew = $run3Start + ($newRun2 + $newRun3).Length
$splitPoint = $run3Start + $newRun2.Length

$d.Bookmarks.Add("ZZ_SPLIT_N", $d.Range($splitPoint, $splitPoint)) | Out-Null

$rngA = $d.Range($run3Start, $splitPoint)
$rngA.Text = $newRun2

$bmN = $d.Bookmarks("ZZ_SPLIT_N")
$rngB = $d.Range($bmN.Start, $run3EndNew)
$rngB.Text = $newRun3

# 5) Clean up the helper bookmarks plus the original hidden "_GoBack"
#    bookmark that used to sit at the old split point.
$d.Bookmarks("ZZ_SPLIT_M").Delete()
$d.Bookmarks("ZZ_SPLIT_N").Delete()
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # already gone / never existed - nothing to do
}

$checkFull = $d.Content.Text
$checkIdx = $checkFull.IndexOf("If you get any")
Write-Output "Updated paragraph text:"
Write-Output $d.Range($checkIdx, $checkFull.Length).Text
